# Generate and save output file after processing
#
# This script reproduces the following structural change to Sheet1:
#  - Three new columns are inserted at R:T (general_college_subjects.history,
#    general_college_subjects.electives, general_college_subjects.cs), which
#    pushes the former R:AE columns right to U:AH.
#  - The new R2/S2/T2 data cells get numeric values 2, 1, 0.
#  - The string value "Unknown" in D2:J2 is normalized to lowercase "unknown".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank columns before the existing "general_college_subjects.arts"
# column (R), shifting R:AE -> U:AH and preserving their values/styles.
$ws.Columns("R:T").Insert()

# Populate the headers for the three newly inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# Populate the data row for the three newly inserted columns.
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 0

# Normalize the "Unknown" placeholders to lowercase "unknown".
$ws.Range("D2:J2").Value = "unknown"

$wb.Save()
